$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N (14), shifting existing
# N:P data to O:Q. This matches the native Excel "Insert Column" gesture.
$ws.Columns.Item(14).Insert()

# The newly inserted column inherits the width of the column to its
# left (M, width 10.7109375) instead of the old column N's width.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet/tab, with R7 selected.
$ws.Activate()
$null = $ws.Range("R7").Select()
